$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "AutoGPT 알아보기"
$ws.Range("E7").Value = "https://jayhey.github.io/deep%20learning/2023/05/07/autogpt/"

$ws.Range("D9").Value = "[스타트업] VC교육으로 얻을 성과"
$ws.Range("E9").Value = "https://pdsi.pabii.com/growing-startup-supporting-vc-10w-3/#utm_source=rss&utm_medium=rss&utm_campaign=growing-startup-supporting-vc-10w-3"

$ws.Range("D32").Value = "Pandas-AI (pandas 활용을 chatGPT 명령에 따라 실행)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/429"
